# 20240703 DataSet.AddBehavior: HitMiss support table input
# ReadQueryTable: support merged cells
# TrialwiseEventPlot: ExcludedEvents -> name-value parameter, add ExcludedTrials
#
# This script reshapes the "query table" sheet so that repeated values in
# columns A-E are shown once and the corresponding cells are merged
# (center / middle-center aligned), matching how ReadQueryTable now
# understands merged-cell input.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# Each entry is an A1 range whose cells should end up merged into one,
# centered both horizontally and vertically. The top-left cell of each
# range keeps its existing value; Merge() clears the other cells in the
# range automatically (matching native Excel behavior).
$mergeRanges = @(
    "B2:B4",
    "C2:C3",
    "D2:D3",
    "C4:C5",
    "A6:A7",
    "C6:C7",
    "E5:E7",
    "A8:A10",
    "B7:B9",
    "C9:C10",
    "E8:E10",
    "A11:A12",
    "B11:B14",
    "C12:C13",
    "A13:A15",
    "C14:C15",
    "E11:E15"
)

foreach ($addr in $mergeRanges) {
    $r = $ws.Range($addr)
    $r.HorizontalAlignment = $xlCenter
    $r.VerticalAlignment = $xlCenter
    $r.Merge()
}

# Restore the selection to match the saved workbook state.
$ws.Range("E19").Select()
